# Fund sheet (基金受益憑證, sheet index 5): add proper header row labels and
# append the property_category/category/date/legislator_name/legislator_id/
# source_file/index metadata columns (I:O), matching the pattern already
# used on the stock (股票) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1): replace the old (buggy, data-duplicate) header
# cells with the real column-name headers, and extend with the new columns.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"

# New header cells I1:O1 need to pick up the same bold/bordered style (style
# index 1) used by the rest of row 1 -- copy formatting from an existing
# header cell before writing the value.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("I1:O1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# --- Data rows (2-11): fix up the dealer name (column D held the fund name
# previously; the real dealer/bank name now goes there) and append the
# metadata columns I:O, copying their style (index 2) from an existing data
# cell on the same row first.
$dealer = "台北富邦銀行"
$fundNames = @{
    2  = "貝萊德世界礦業基金"
    3  = "天達動力資源基金"
    4  = "施羅德環球美元流動基金"
    5  = "貝萊德世界礦業基金"
    6  = "瑞銀美元基金"
    7  = "貝萊德美元儲備基金"
    8  = "貝萊德環球資產配置基金"
    9  = "景順美元儲備基"
    10 = "貝萊德新興市場基金"
    11 = "貝萊德世界礦業基金"
}

for ($r = 2; $r -le 11; $r++) {
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("I$r:O$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("B$r").Value = $fundNames[$r]
    $ws.Range("D$r").Value = $dealer

    $idx = $ws.Range("A$r").Value

    $ws.Range("I$r").Value = "fund"
    $ws.Range("J$r").Value = "normal"
    $ws.Range("K$r").Value = "2011-11-22"
    $ws.Range("L$r").Value = "徐少萍"
    $ws.Range("M$r").Value = 726
    $ws.Range("N$r").Value = "tmpc12c1"
    $ws.Range("O$r").Value = $idx
}
